# Refatoracao do exercicio 17 (Teste de Mesa)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Insert a new row at 81 (pushes the old row 81 down to row 82)
#    so the exercise-17 table grows from 2 data rows to 3.
# ---------------------------------------------------------------
$ws.Rows.Item(81).Insert()

# ---------------------------------------------------------------
# 2) Header row (row 79): drop "fator", "api" moves to column F,
#    a brand new "pontosExtras" column takes the old "api" slot (G).
# ---------------------------------------------------------------
$ws.Range("F79").Value = "api"
$ws.Range("G79").Value = "pontosExtras"

# ---------------------------------------------------------------
# 3) Give the new row 81 the same border formatting as the rest of
#    the table by copying the format from row 80.
# ---------------------------------------------------------------
$ws.Range("A80:I80").Copy()
$ws.Range("A81:I81").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 4) Data values.
#    Row 80 (p1=8 e1=6 e2=8 media=3.8 max=1 api=8 pontosExtras=0.5 sub=0 mediaFinal=8.3)
# ---------------------------------------------------------------
$ws.Range("E80").Value = 1
$ws.Range("F80").Value = 8
$ws.Range("I80").Value = "8.3"

# Row 81 (new): p1=8 e1=6 e2=7 media=3.7 max=1 api=3 pontosExtras=0.5 sub=5 mediaFinal=7.2
$ws.Range("A81").Value = 8
$ws.Range("B81").Value = 6
$ws.Range("C81").Value = 7
$ws.Range("D81").Value = "3.7"
$ws.Range("E81").Value = 1
$ws.Range("F81").Value = 3
$ws.Range("H81").Value = 5
$ws.Range("I81").Value = "7.2"

# Row 82 (was row 81 before the insert): p1=5 e1=6 e2=4 media=2.5 max=0 api=0 pontosExtras=0.5 sub=8 mediaFinal=5.4
$ws.Range("I82").Value = "5.4"

# ---------------------------------------------------------------
# 5) "pontosExtras" column (G) data cells hold the text "0.5" but
#    are displayed with a 0.00 number format, right aligned.
# ---------------------------------------------------------------
$ws.Range("G80").NumberFormat = "@"
$ws.Range("G80").Value = "0.5"
$ws.Range("G81").NumberFormat = "@"
$ws.Range("G81").Value = "0.5"
$ws.Range("G82").NumberFormat = "@"
$ws.Range("G82").Value = "0.5"
$ws.Range("G80:G82").NumberFormat = "0.00"

# ---------------------------------------------------------------
# 6) Whole table (header + data, A79:I82) is right aligned now.
# ---------------------------------------------------------------
$ws.Range("A79:I82").HorizontalAlignment = -4152

# ---------------------------------------------------------------
# 7) Column widths: H grows, new column I keeps H's old width.
# ---------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 11.92
$ws.Columns.Item(9).ColumnWidth = 9.17

# ---------------------------------------------------------------
# 8) View state: scroll down a bit, select G84.
# ---------------------------------------------------------------
$ws.Range("G84").Select()

Write-Host "edit complete"
